# Fix double entries for 2023-05-15 (serial 45062) in the BL sheet.
# Rows 479, 481, 483, ... 511 are exact duplicates of the row directly
# above them (478, 480, 482, ... 510). Delete the duplicate (second)
# occurrence of each pair, working from the bottom up so row numbers of
# not-yet-processed rows are unaffected by the deletions above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = 511, 509, 507, 505, 503, 501, 499, 497, 495, 493, 491, 489, 487, 485, 483, 481, 479

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
